# Apply the edit described by the diff:
#  1. For every 4-row year-group (A,B,C,D) in rows 2..69, swap the data of
#     the "B" and "C" sub-rows (columns A:E), leaving the "A" and "D"
#     sub-rows untouched.
#  2. Remove columns F and G entirely (产销率 / 销售量 duplicate columns),
#     which also shrinks the sheet dimension from A1:G69 to A1:E69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 69
$lastCol = 5   # columns A..E hold the data that must survive

for ($r = 2; $r -le $lastRow; $r += 4) {
    $rowB = $r + 1
    $rowC = $r + 2

    if ($rowC -gt $lastRow) {
        continue
    }

    # Stash row B's A:E values, then overwrite B with C's values, then
    # write the stashed B values into C.
    for ($c = 1; $c -le $lastCol; $c++) {
        $bVal = $ws.Cells.Item($rowB, $c).Value2
        $cVal = $ws.Cells.Item($rowC, $c).Value2

        $ws.Cells.Item($rowB, $c).Value2 = $cVal
        $ws.Cells.Item($rowC, $c).Value2 = $bVal
    }
}

# Drop the now-redundant F/G columns (dimension becomes A1:E69).
$ws.Range("F1:G$lastRow").EntireColumn.Delete()
